# Scheduled Gungnir Profits data runner: refresh cached market price /
# profit figures (currentAveragePrice*, LevePrice*, LeveProfit*) per leve row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3066.6667
$ws.Range("I64").Value = 2866.6667
$ws.Range("J64").Value = 3466.6667
$ws.Range("K64").Value = 2866.6667
$ws.Range("L64").Value = 3466.6667
$ws.Range("M64").Value = -2618.6667
$ws.Range("N64").Value = -3962.6667
$ws.Range("H67").Value = 3066.6667
$ws.Range("I67").Value = 2866.6667
$ws.Range("J67").Value = 3466.6667
$ws.Range("K67").Value = 2866.6667
$ws.Range("L67").Value = 3466.6667
$ws.Range("M67").Value = -2008.6667
$ws.Range("N67").Value = -5182.6667
$ws.Range("H96").Value = 965
$ws.Range("I96").Value = 589.5714
$ws.Range("J96").Value = 1403
$ws.Range("K96").Value = 1768.7142
$ws.Range("L96").Value = 4209
$ws.Range("M96").Value = -395.7142000000001
$ws.Range("N96").Value = -6955
$ws.Range("H141").Value = 4700.0835
$ws.Range("I141").Value = 3274
$ws.Range("J141").Value = 5413.125
$ws.Range("K141").Value = 9822
$ws.Range("L141").Value = 16239.375
$ws.Range("M141").Value = -4642
$ws.Range("N141").Value = -26599.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 817.35297
$ws.Range("I2").Value = 782.8570999999999
$ws.Range("J2").Value = 841.5
$ws.Range("K2").Value = 782.8570999999999
$ws.Range("L2").Value = 841.5
$ws.Range("M2").Value = -669.8570999999999
$ws.Range("N2").Value = -1067.5
$ws.Range("H45").Value = 56709.832
$ws.Range("I45").Value = 100830
$ws.Range("J45").Value = 1559.625
$ws.Range("K45").Value = 100830
$ws.Range("L45").Value = 1559.625
$ws.Range("M45").Value = -100453
$ws.Range("N45").Value = -2313.625
$ws.Range("H64").Value = 17786.334
$ws.Range("J64").Value = 17786.334
$ws.Range("L64").Value = 17786.334
$ws.Range("N64").Value = -18282.334
$ws.Range("H67").Value = 17786.334
$ws.Range("J67").Value = 17786.334
$ws.Range("L67").Value = 17786.334
$ws.Range("N67").Value = -19502.334
$ws.Range("H116").Value = 817.35297
$ws.Range("I116").Value = 782.8570999999999
$ws.Range("J116").Value = 841.5
$ws.Range("K116").Value = 782.8570999999999
$ws.Range("L116").Value = 841.5
$ws.Range("M116").Value = 1511.1429
$ws.Range("N116").Value = -5429.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 817.35297
$ws.Range("I3").Value = 782.8570999999999
$ws.Range("J3").Value = 841.5
$ws.Range("K3").Value = 782.8570999999999
$ws.Range("L3").Value = 841.5
$ws.Range("M3").Value = -668.8570999999999
$ws.Range("N3").Value = -1069.5
$ws.Range("H62").Value = 19800
$ws.Range("J62").Value = 19800
$ws.Range("L62").Value = 19800
$ws.Range("N62").Value = -21172
$ws.Range("H65").Value = 19800
$ws.Range("J65").Value = 19800
$ws.Range("L65").Value = 59400
$ws.Range("N65").Value = -66264

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1613.7894
$ws.Range("I16").Value = 1105.4445
$ws.Range("J16").Value = 2071.3
$ws.Range("K16").Value = 1105.4445
$ws.Range("L16").Value = 2071.3
$ws.Range("M16").Value = -818.4445000000001
$ws.Range("N16").Value = -2645.3
$ws.Range("H62").Value = 3734.875
$ws.Range("I62").Value = 2849.8333
$ws.Range("J62").Value = 6390
$ws.Range("K62").Value = 2849.8333
$ws.Range("L62").Value = 6390
$ws.Range("M62").Value = -2225.8333
$ws.Range("N62").Value = -7638
$ws.Range("H65").Value = 3734.875
$ws.Range("I65").Value = 2849.8333
$ws.Range("J65").Value = 6390
$ws.Range("K65").Value = 14249.1665
$ws.Range("L65").Value = 31950
$ws.Range("M65").Value = -11129.1665
$ws.Range("N65").Value = -38190
$ws.Range("H113").Value = 1613.7894
$ws.Range("I113").Value = 1105.4445
$ws.Range("J113").Value = 2071.3
$ws.Range("K113").Value = 1105.4445
$ws.Range("L113").Value = 2071.3
$ws.Range("M113").Value = 1064.5555
$ws.Range("N113").Value = -6411.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 30308200
$ws.Range("I5").Value = 44444870
$ws.Range("J5").Value = 15328.571
$ws.Range("K5").Value = 133334610
$ws.Range("L5").Value = 45985.713
$ws.Range("M5").Value = -133334498
$ws.Range("N5").Value = -46209.713
$ws.Range("H107").Value = 33336840
$ws.Range("I107").Value = 219.64285
$ws.Range("J107").Value = 59265324
$ws.Range("K107").Value = 658.9285500000001
$ws.Range("L107").Value = 177795972
$ws.Range("M107").Value = 1261.07145
$ws.Range("N107").Value = -177799812
$ws.Range("H122").Value = 13592001
$ws.Range("I122").Value = 89285960
$ws.Range("J122").Value = 5906.359
$ws.Range("K122").Value = 803573640
$ws.Range("L122").Value = 53157.231
$ws.Range("M122").Value = -803571190
$ws.Range("N122").Value = -58057.231
$ws.Range("H135").Value = 30308200
$ws.Range("I135").Value = 44444870
$ws.Range("J135").Value = 15328.571
$ws.Range("K135").Value = 400003830
$ws.Range("L135").Value = 137957.139
$ws.Range("M135").Value = -400001295
$ws.Range("N135").Value = -143027.139
$ws.Range("H136").Value = 13161396
$ws.Range("I136").Value = 27779592
$ws.Range("K136").Value = 83338776
$ws.Range("M136").Value = -83333676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1140.2142
$ws.Range("I113").Value = 1028.5714
$ws.Range("J113").Value = 1251.8572
$ws.Range("K113").Value = 1028.5714
$ws.Range("L113").Value = 1251.8572
$ws.Range("M113").Value = 1141.4286
$ws.Range("N113").Value = -5591.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1496.2307
$ws.Range("I68").Value = 1496.75
$ws.Range("J68").Value = 1490
$ws.Range("K68").Value = 1496.75
$ws.Range("L68").Value = 1490
$ws.Range("M68").Value = -747.75
$ws.Range("N68").Value = -2988
$ws.Range("H71").Value = 1496.2307
$ws.Range("I71").Value = 1496.75
$ws.Range("J71").Value = 1490
$ws.Range("K71").Value = 7483.75
$ws.Range("L71").Value = 7450
$ws.Range("M71").Value = -3739.75
$ws.Range("N71").Value = -14938

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4077.7778
$ws.Range("I62").Value = 4283.3335
$ws.Range("J62").Value = 3666.6667
$ws.Range("K62").Value = 4283.3335
$ws.Range("L62").Value = 3666.6667
$ws.Range("M62").Value = -3659.3335
$ws.Range("N62").Value = -4914.6667
$ws.Range("H65").Value = 4077.7778
$ws.Range("I65").Value = 4283.3335
$ws.Range("J65").Value = 3666.6667
$ws.Range("K65").Value = 21416.6675
$ws.Range("L65").Value = 18333.3335
$ws.Range("M65").Value = -18296.6675
$ws.Range("N65").Value = -24573.3335

